$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "comptence_ok" variable row (row 17), matching the other
# rows in the "Oui" (displayed variable) block above it.
$ws.Range("A17").Value = "Oui"
$ws.Range("B17").Value = "comptence_ok"
$ws.Range("C17").Value = "La proportion de jeunes estimant être employés à leur niveau de compétence"
